$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values stored as text. Prefixing with a
# leading apostrophe forces Excel to keep them as text (matching the
# original inline-string cell type) instead of coercing them to numbers.

$ws.Range("D2").Value = "'276.91"
$ws.Range("D3").Value = "'21.16"
$ws.Range("D4").Value = "'6.270"
$ws.Range("D5").Value = "'0.06209"
$ws.Range("D6").Value = "'3.555"
$ws.Range("D7").Value = "'1.533"
$ws.Range("D8").Value = "'6.580"
$ws.Range("D9").Value = "'0.8280"
$ws.Range("D10").Value = "'0.1667"
$ws.Range("D11").Value = "'0.08291"
$ws.Range("D12").Value = "'0.03515"
$ws.Range("D13").Value = "'0.03167"
$ws.Range("D14").Value = "'0.09181"
$ws.Range("D15").Value = "'3.766"
$ws.Range("D16").Value = "'0.001628"
$ws.Range("D17").Value = "'0.04686"
$ws.Range("D18").Value = "'0.006395"
$ws.Range("D19").Value = "'0.006215"
$ws.Range("D20").Value = "'0.001068"
$ws.Range("D21").Value = "'0.0001500"
$ws.Range("D22").Value = "'3.721"
$ws.Range("D23").Value = "'2.313"
$ws.Range("D24").Value = "'0.01397"
$ws.Range("D28").Value = "'0.0002736"
$ws.Range("D40").Value = "'0.04745"

# Rows 41-43 rotate their Coin/Link/Price/Volume data:
#   new row41 <- old row42 data (KickToken), with an updated price
#   new row42 <- old row43 data (BKEXToken), with an updated price
#   new row43 <- old row41 data (CEJI), with an updated price
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007058"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1122"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003519"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "'0.01136"
$ws.Range("D45").Value = "'0.00006354"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D48").Value = "'0.9901"
$ws.Range("D51").Value = "'0.01240"
